$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the invoice number value
$ws.Range("A2").Value = "INV-00010"

# Set column A width (~21.66 chars; engine quantizes ColumnWidth to 1/6-char
# pixel steps, so 20.8333... is the input that round-trips closest to the
# target stored width of 21.6640625)
$ws.Columns.Item(1).ColumnWidth = 20.8333333333333

# Update selection
$ws.Range("C13").Select()
